$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug correction in investigator name fields (row 2): Xing H. Rosie -> Robert Jones
$ws.Range("D2").Value = "Robert Jones"
$ws.Range("G2").Value = "Robert"
$ws.Range("I2").Value = "Jones"

# Update view: scroll to G1, select I5
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I5").Select()
